$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ordered?" column (A) for the Component Order BOM table (rows 57-80)
$ws.Range("A57").Value = "Ordered?"
$ws.Range("A58").Value = "n/a"
$ws.Range("A59").Value = "Y"
$ws.Range("A60").Value = "no"
$ws.Range("A61").Value = "Y"
$ws.Range("A62").Value = "Y"
$ws.Range("A63").Value = "Y"
$ws.Range("A64").Value = "Y"
$ws.Range("A65").Value = "Y"
$ws.Range("A66").Value = "no"
$ws.Range("A67").Value = "Y"
$ws.Range("A68").Value = "Y"
$ws.Range("A69").Value = "Y"
$ws.Range("A70").Value = "Y"
$ws.Range("A71").Value = "no"
$ws.Range("A72").Value = "Y"
$ws.Range("A73").Value = "Y"
$ws.Range("A74").Value = "Y"
$ws.Range("A75").Value = "Y"
$ws.Range("A76").Value = "Y"
$ws.Range("A77").Value = "no"
$ws.Range("A78").Value = "Y"
$ws.Range("A79").Value = "Y"
$ws.Range("A80").Value = "Y"

# New note row below the totals, with layout-fix instructions (commit: "Fix DCDC track layout")
$ws.Range("B84").Value = "Move capacitor nearer the dcdc chip, keep inductor feedback as short as possible, don't route under the inductor, connect en pin to via and leave floating"

# Leave the cursor where the user finished typing, matching the saved selection
$ws.Range("B84").Select()
